$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Update the "historic" period dates in row 8 ---
$ws.Range("B8").Value = 44743   # Fecha de inicio del periodo que se informa
$ws.Range("C8").Value = 44834   # Fecha de termino del periodo que se informa
$ws.Range("I8").Value = 44844   # Fecha de validacion
$ws.Range("J8").Value = 44844   # Fecha de actualizacion

# --- Widen column K (11) to match the new content width ---
$ws.Columns.Item(11).ColumnWidth = 48.83

# --- Update the view: scrolled position / active selection ---
$ws.Activate()
$ws.Range("C14").Select()

$wb.Save()
